$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Step 1: stash format templates from pristine rows into scratch area ----
$ws.Range("A3:E3").Copy() | Out-Null
$ws.Range("AA1:AE1").PasteSpecial(-4122) | Out-Null
$ws.Range("A9:E9").Copy() | Out-Null
$ws.Range("AA2:AE2").PasteSpecial(-4122) | Out-Null
$ws.Range("A4:D4").Copy() | Out-Null
$ws.Range("AA3:AD3").PasteSpecial(-4122) | Out-Null
$ws.Range("C4").Copy() | Out-Null
$ws.Range("AE3").PasteSpecial(-4122) | Out-Null
$ws.Range("A6:E6").Copy() | Out-Null
$ws.Range("AA4:AE4").PasteSpecial(-4122) | Out-Null
$ws.Range("A7:E7").Copy() | Out-Null
$ws.Range("AA5:AE5").PasteSpecial(-4122) | Out-Null
$ws.Range("A8:E8").Copy() | Out-Null
$ws.Range("AA6:AE6").PasteSpecial(-4122) | Out-Null
$ws.Range("G22:H22").Copy() | Out-Null
$ws.Range("AA7:AB7").PasteSpecial(-4122) | Out-Null
$ws.Range("G23:H23").Copy() | Out-Null
$ws.Range("AA8:AB8").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# ---- Step 2: unmerge existing merged cells ----
$ws.Range("B6:B8").UnMerge() | Out-Null
$ws.Range("B10:B12").UnMerge() | Out-Null
$ws.Range("B13:B15").UnMerge() | Out-Null
$ws.Range("B16:B18").UnMerge() | Out-Null
$ws.Range("B19:B20").UnMerge() | Out-Null

# ---- Step 3: clear old content area (A3:H27) ----
$ws.Range("A3:H27").Clear() | Out-Null

# ---- Step 4: rebuild header row 3 ----
$ws.Range("AA1:AE1").Copy() | Out-Null
$ws.Range("A3:E3").PasteSpecial(-4122) | Out-Null
$ws.Range("B3").Value2 = "Piso"
$ws.Range("C3").Value2 = "Sector"
$ws.Range("D3").Value2 = "Tipo Matafuego"
$ws.Range("E3").Value2 = "Cantidad"

# ---- Step 5: rebuild data rows 4-25 ----
# Row 4
$ws.Range("AB2").Copy() | Out-Null
$ws.Range("B4").PasteSpecial(-4122) | Out-Null
$ws.Range("B4").Value2 = "Edificio Expedición"
$ws.Range("AC4:AD4").Copy() | Out-Null
$ws.Range("C4:D4").PasteSpecial(-4122) | Out-Null
$ws.Range("AE4").Copy() | Out-Null
$ws.Range("E4").PasteSpecial(-4122) | Out-Null
$ws.Range("C4").Value2 = "Expedición"
$ws.Range("D4").Value2 = "Clase ABC"
$ws.Range("E4").Value2 = 5

# Row 5
$ws.Range("AA3").Copy() | Out-Null
$ws.Range("A5").PasteSpecial(-4122) | Out-Null
$ws.Range("AB3").Copy() | Out-Null
$ws.Range("B5").PasteSpecial(-4122) | Out-Null
$ws.Range("B5").Value2 = "Primer Subsuelo"
$ws.Range("AC3:AD3").Copy() | Out-Null
$ws.Range("C5:D5").PasteSpecial(-4122) | Out-Null
$ws.Range("AE3").Copy() | Out-Null
$ws.Range("E5").PasteSpecial(-4122) | Out-Null
$ws.Range("C5").Value2 = "Estacionamiento"
$ws.Range("D5").Value2 = "Clase ABC"
$ws.Range("E5").Value2 = 1

# Row 6
$ws.Range("AA3").Copy() | Out-Null
$ws.Range("A6").PasteSpecial(-4122) | Out-Null
$ws.Range("AB3").Copy() | Out-Null
$ws.Range("B6").PasteSpecial(-4122) | Out-Null
$ws.Range("B6").Value2 = "Segundo Subsuelo"
$ws.Range("AC3:AD3").Copy() | Out-Null
$ws.Range("C6:D6").PasteSpecial(-4122) | Out-Null
$ws.Range("AE3").Copy() | Out-Null
$ws.Range("E6").PasteSpecial(-4122) | Out-Null
$ws.Range("C6").Value2 = "Estacionamiento"
$ws.Range("D6").Value2 = "Clase ABC"
$ws.Range("E6").Value2 = 1

# Row 7
$ws.Range("AB4").Copy() | Out-Null
$ws.Range("B7").PasteSpecial(-4122) | Out-Null
$ws.Range("B7").Value2 = "Planta Baja "
$ws.Range("AC4:AD4").Copy() | Out-Null
$ws.Range("C7:D7").PasteSpecial(-4122) | Out-Null
$ws.Range("AE4").Copy() | Out-Null
$ws.Range("E7").PasteSpecial(-4122) | Out-Null
$ws.Range("C7").Value2 = "Buffet"
$ws.Range("D7").Value2 = "Clase K, Clase AC"
$ws.Range("E7").Value2 = 2

# Row 8
$ws.Range("AB5").Copy() | Out-Null
$ws.Range("B8").PasteSpecial(-4122) | Out-Null
$ws.Range("AC4:AD4").Copy() | Out-Null
$ws.Range("C8:D8").PasteSpecial(-4122) | Out-Null
$ws.Range("AE4").Copy() | Out-Null
$ws.Range("E8").PasteSpecial(-4122) | Out-Null
$ws.Range("C8").Value2 = "Servidores"
$ws.Range("D8").Value2 = "Clase AC"
$ws.Range("E8").Value2 = 1

# Row 9
$ws.Range("AB6").Copy() | Out-Null
$ws.Range("B9").PasteSpecial(-4122) | Out-Null
$ws.Range("AC4:AD4").Copy() | Out-Null
$ws.Range("C9:D9").PasteSpecial(-4122) | Out-Null
$ws.Range("AE4").Copy() | Out-Null
$ws.Range("E9").PasteSpecial(-4122) | Out-Null
$ws.Range("C9").Value2 = "Mantenimiento"
$ws.Range("D9").Value2 = "Clase ABC"
$ws.Range("E9").Value2 = 1

# Row 10
$ws.Range("AB4").Copy() | Out-Null
$ws.Range("B10").PasteSpecial(-4122) | Out-Null
$ws.Range("B10").Value2 = "Primer Piso"
$ws.Range("AC4:AD4").Copy() | Out-Null
$ws.Range("C10:D10").PasteSpecial(-4122) | Out-Null
$ws.Range("AE4").Copy() | Out-Null
$ws.Range("E10").PasteSpecial(-4122) | Out-Null
$ws.Range("C10").Value2 = "Producto"
$ws.Range("D10").Value2 = "Clase AC"
$ws.Range("E10").Value2 = 1

# Row 11
$ws.Range("AB5").Copy() | Out-Null
$ws.Range("B11").PasteSpecial(-4122) | Out-Null
$ws.Range("AC4:AD4").Copy() | Out-Null
$ws.Range("C11:D11").PasteSpecial(-4122) | Out-Null
$ws.Range("AE4").Copy() | Out-Null
$ws.Range("E11").PasteSpecial(-4122) | Out-Null
$ws.Range("C11").Value2 = "Redacción"
$ws.Range("D11").Value2 = "Clase AC"
$ws.Range("E11").Value2 = 1

# Row 12
$ws.Range("AB6").Copy() | Out-Null
$ws.Range("B12").PasteSpecial(-4122) | Out-Null
$ws.Range("AC4:AD4").Copy() | Out-Null
$ws.Range("C12:D12").PasteSpecial(-4122) | Out-Null
$ws.Range("AE4").Copy() | Out-Null
$ws.Range("E12").PasteSpecial(-4122) | Out-Null
$ws.Range("C12").Value2 = "kitchenette"
$ws.Range("D12").Value2 = "Clase K"
$ws.Range("E12").Value2 = 1

# Row 13
$ws.Range("AB4").Copy() | Out-Null
$ws.Range("B13").PasteSpecial(-4122) | Out-Null
$ws.Range("B13").Value2 = "Segundo Piso"
$ws.Range("AC4:AD4").Copy() | Out-Null
$ws.Range("C13:D13").PasteSpecial(-4122) | Out-Null
$ws.Range("AE4").Copy() | Out-Null
$ws.Range("E13").PasteSpecial(-4122) | Out-Null
$ws.Range("C13").Value2 = "Área Digital"
$ws.Range("D13").Value2 = "Clase AC"
$ws.Range("E13").Value2 = 2

# Row 14
$ws.Range("AB5").Copy() | Out-Null
$ws.Range("B14").PasteSpecial(-4122) | Out-Null
$ws.Range("AC4:AD4").Copy() | Out-Null
$ws.Range("C14:D14").PasteSpecial(-4122) | Out-Null
$ws.Range("AE4").Copy() | Out-Null
$ws.Range("E14").PasteSpecial(-4122) | Out-Null
$ws.Range("C14").Value2 = "Pasillo Sala y Oficinas"
$ws.Range("D14").Value2 = "Clase AC"
$ws.Range("E14").Value2 = 1

# Row 15
$ws.Range("AB5").Copy() | Out-Null
$ws.Range("B15").PasteSpecial(-4122) | Out-Null
$ws.Range("AC4:AD4").Copy() | Out-Null
$ws.Range("C15:D15").PasteSpecial(-4122) | Out-Null
$ws.Range("AE4").Copy() | Out-Null
$ws.Range("E15").PasteSpecial(-4122) | Out-Null
$ws.Range("C15").Value2 = "kitchenette"
$ws.Range("D15").Value2 = "Clase K"
$ws.Range("E15").Value2 = 1

# Row 16
$ws.Range("AB6").Copy() | Out-Null
$ws.Range("B16").PasteSpecial(-4122) | Out-Null
$ws.Range("AC4:AD4").Copy() | Out-Null
$ws.Range("C16:D16").PasteSpecial(-4122) | Out-Null
$ws.Range("AE4").Copy() | Out-Null
$ws.Range("E16").PasteSpecial(-4122) | Out-Null
$ws.Range("C16").Value2 = "Escaleras"
$ws.Range("D16").Value2 = "Clase AC"
$ws.Range("E16").Value2 = 1

# Row 17
$ws.Range("AB4").Copy() | Out-Null
$ws.Range("B17").PasteSpecial(-4122) | Out-Null
$ws.Range("B17").Value2 = "Tercer Piso"
$ws.Range("AC4:AD4").Copy() | Out-Null
$ws.Range("C17:D17").PasteSpecial(-4122) | Out-Null
$ws.Range("AE4").Copy() | Out-Null
$ws.Range("E17").PasteSpecial(-4122) | Out-Null
$ws.Range("C17").Value2 = "Capacitación"
$ws.Range("D17").Value2 = "Clase AC"
$ws.Range("E17").Value2 = 2

# Row 18
$ws.Range("AB5").Copy() | Out-Null
$ws.Range("B18").PasteSpecial(-4122) | Out-Null
$ws.Range("AC4:AD4").Copy() | Out-Null
$ws.Range("C18:D18").PasteSpecial(-4122) | Out-Null
$ws.Range("AE4").Copy() | Out-Null
$ws.Range("E18").PasteSpecial(-4122) | Out-Null
$ws.Range("C18").Value2 = "Entrevistas"
$ws.Range("D18").Value2 = "Clase AC"
$ws.Range("E18").Value2 = 1

# Row 19
$ws.Range("AB6").Copy() | Out-Null
$ws.Range("B19").PasteSpecial(-4122) | Out-Null
$ws.Range("AC4:AD4").Copy() | Out-Null
$ws.Range("C19:D19").PasteSpecial(-4122) | Out-Null
$ws.Range("AE4").Copy() | Out-Null
$ws.Range("E19").PasteSpecial(-4122) | Out-Null
$ws.Range("C19").Value2 = "kitchenette"
$ws.Range("D19").Value2 = "Clase K"
$ws.Range("E19").Value2 = 1

# Row 20
$ws.Range("AB4").Copy() | Out-Null
$ws.Range("B20").PasteSpecial(-4122) | Out-Null
$ws.Range("B20").Value2 = "Cuarto Piso"
$ws.Range("AC4:AD4").Copy() | Out-Null
$ws.Range("C20:D20").PasteSpecial(-4122) | Out-Null
$ws.Range("AE4").Copy() | Out-Null
$ws.Range("E20").PasteSpecial(-4122) | Out-Null
$ws.Range("C20").Value2 = "Sistemas"
$ws.Range("D20").Value2 = "Clase AC"
$ws.Range("E20").Value2 = 1

# Row 21
$ws.Range("AB5").Copy() | Out-Null
$ws.Range("B21").PasteSpecial(-4122) | Out-Null
$ws.Range("AC4:AD4").Copy() | Out-Null
$ws.Range("C21:D21").PasteSpecial(-4122) | Out-Null
$ws.Range("AE4").Copy() | Out-Null
$ws.Range("E21").PasteSpecial(-4122) | Out-Null
$ws.Range("C21").Value2 = "kitchenette"
$ws.Range("D21").Value2 = "Clase K"
$ws.Range("E21").Value2 = 1

# Row 22
$ws.Range("AB6").Copy() | Out-Null
$ws.Range("B22").PasteSpecial(-4122) | Out-Null
$ws.Range("AC4:AD4").Copy() | Out-Null
$ws.Range("C22:D22").PasteSpecial(-4122) | Out-Null
$ws.Range("AE4").Copy() | Out-Null
$ws.Range("E22").PasteSpecial(-4122) | Out-Null
$ws.Range("C22").Value2 = "Escaleras"
$ws.Range("D22").Value2 = "Clase AC"
$ws.Range("E22").Value2 = 1

# Row 23
$ws.Range("AB4").Copy() | Out-Null
$ws.Range("B23").PasteSpecial(-4122) | Out-Null
$ws.Range("B23").Value2 = "Quinto Piso"
$ws.Range("AC4:AD4").Copy() | Out-Null
$ws.Range("C23:D23").PasteSpecial(-4122) | Out-Null
$ws.Range("AE4").Copy() | Out-Null
$ws.Range("E23").PasteSpecial(-4122) | Out-Null
$ws.Range("C23").Value2 = "Gerencia General"
$ws.Range("D23").Value2 = "Clase AC"
$ws.Range("E23").Value2 = 1

# Row 24
$ws.Range("AB5").Copy() | Out-Null
$ws.Range("B24").PasteSpecial(-4122) | Out-Null
$ws.Range("AC4:AD4").Copy() | Out-Null
$ws.Range("C24:D24").PasteSpecial(-4122) | Out-Null
$ws.Range("AE4").Copy() | Out-Null
$ws.Range("E24").PasteSpecial(-4122) | Out-Null
$ws.Range("C24").Value2 = "kitchenette"
$ws.Range("D24").Value2 = "Clase K"
$ws.Range("E24").Value2 = 1

# Row 25
$ws.Range("AB6").Copy() | Out-Null
$ws.Range("B25").PasteSpecial(-4122) | Out-Null
$ws.Range("AC4:AD4").Copy() | Out-Null
$ws.Range("C25:D25").PasteSpecial(-4122) | Out-Null
$ws.Range("AE4").Copy() | Out-Null
$ws.Range("E25").PasteSpecial(-4122) | Out-Null
$ws.Range("C25").Value2 = "Sala de Reuniones"
$ws.Range("D25").Value2 = "Clase AC"
$ws.Range("E25").Value2 = 1

# ---- Step 6: re-merge B column groups ----
$ws.Range("B7:B9").Merge() | Out-Null
$ws.Range("B10:B12").Merge() | Out-Null
$ws.Range("B13:B16").Merge() | Out-Null
$ws.Range("B17:B19").Merge() | Out-Null
$ws.Range("B20:B22").Merge() | Out-Null
$ws.Range("B23:B25").Merge() | Out-Null

# ---- Step 7: rebuild side table (Utilizados para), rows 25-30 ----
$ws.Range("AA7:AB7").Copy() | Out-Null
$ws.Range("G25:H25").PasteSpecial(-4122) | Out-Null
$ws.Range("G25").Value2 = "Tipo Matafuego"
$ws.Range("H25").Value2 = "Utilizados para:"

$ws.Range("AA8:AB8").Copy() | Out-Null
$ws.Range("G26:H26").PasteSpecial(-4122) | Out-Null
$ws.Range("G26").Value2 = "Clase K"
$ws.Range("H26").Value2 = "fuegos de aceites vegetales o grasas animales."

$ws.Range("AA8:AB8").Copy() | Out-Null
$ws.Range("G27:H27").PasteSpecial(-4122) | Out-Null
$ws.Range("G27").Value2 = "Clase A"
$ws.Range("H27").Value2 = "Combustilbes sólidos: Papel, madera, goma"

$ws.Range("AA8:AB8").Copy() | Out-Null
$ws.Range("G28:H28").PasteSpecial(-4122) | Out-Null
$ws.Range("G28").Value2 = "Clase B"
$ws.Range("H28").Value2 = "Combustibles líquidos:pinturas, grasas, solventes, naftas"

$ws.Range("AA8:AB8").Copy() | Out-Null
$ws.Range("G29:H29").PasteSpecial(-4122) | Out-Null
$ws.Range("G29").Value2 = "Clase C"
$ws.Range("H29").Value2 = "Electricidad"

$ws.Range("AA8:AB8").Copy() | Out-Null
$ws.Range("G30:H30").PasteSpecial(-4122) | Out-Null
$ws.Range("G30").Value2 = "Clase D"
$ws.Range("H30").Value2 = "Metales combustibles: magnesio, titanio, zirconio, sodio, potasio"

# ---- Step 8: clean up scratch area ----
$ws.Range("AA1:AE8").Clear() | Out-Null
$excel.CutCopyMode = 0

# ---- Step 9: selection/view ----
$ws.Range("H9").Select() | Out-Null
